# Applies updated transition-matrix values for UMES_A team-specific matrix
# (added more games, sped up simulate game logic, drafted optimization logic)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.1833333333333333
$ws.Range("C2").Value = 0.55
$ws.Range("J2").Value = 0.01666666666666667
$ws.Range("P2").Value = 0.1583333333333333
$ws.Range("S2").Value = 0.09166666666666666
$ws.Range("B3").Value = 0.0145985401459854
$ws.Range("C3").Value = 0.0364963503649635
$ws.Range("J3").Value = 0.04379562043795621
$ws.Range("P3").Value = 0.7153284671532847
$ws.Range("S3").Value = 0.1897810218978102
$ws.Range("J4").Value = 0.04166666666666666
$ws.Range("P4").Value = 0.7083333333333334
$ws.Range("S4").Value = 0.25
$ws.Range("B6").Value = 0.08280254777070063
$ws.Range("D6").Value = 0.01273885350318471
$ws.Range("F6").Value = 0.02547770700636943
$ws.Range("J6").Value = 0.2929936305732484
$ws.Range("O6").Value = 0.02547770700636943
$ws.Range("Q6").Value = 0.1273885350318471
$ws.Range("R6").Value = 0.03821656050955414
$ws.Range("S6").Value = 0.3949044585987261
$ws.Range("B7").Value = 0.1153846153846154
$ws.Range("D7").Value = 0.00641025641025641
$ws.Range("F7").Value = 0.02564102564102564
$ws.Range("J7").Value = 0.141025641025641
$ws.Range("O7").Value = 0.0641025641025641
$ws.Range("Q7").Value = 0.1282051282051282
$ws.Range("R7").Value = 0.07051282051282051
$ws.Range("S7").Value = 0.4487179487179487
$ws.Range("B8").Value = 0.1005434782608696
$ws.Range("D8").Value = 0.008152173913043478
$ws.Range("F8").Value = 0.0516304347826087
$ws.Range("J8").Value = 0.125
$ws.Range("O8").Value = 0.01902173913043478
$ws.Range("Q8").Value = 0.1684782608695652
$ws.Range("R8").Value = 0.07880434782608696
$ws.Range("S8").Value = 0.4483695652173913
$ws.Range("B9").Value = 0.09848484848484848
$ws.Range("D9").Value = 0.007575757575757576
$ws.Range("F9").Value = 0.04545454545454546
$ws.Range("J9").Value = 0.09090909090909091
$ws.Range("O9").Value = 0.03787878787878788
$ws.Range("Q9").Value = 0.1590909090909091
$ws.Range("R9").Value = 0.07575757575757576
$ws.Range("S9").Value = 0.4848484848484849
$ws.Range("B10").Value = 0.09893307468477207
$ws.Range("D10").Value = 0.01648884578079534
$ws.Range("E10").Value = 0.001939864209505335
$ws.Range("F10").Value = 0.05625606207565471
$ws.Range("J10").Value = 0.1251212415130941
$ws.Range("O10").Value = 0.008729388942774006
$ws.Range("Q10").Value = 0.2201745877788555
$ws.Range("R10").Value = 0.06498545101842872
$ws.Range("S10").Value = 0.4073714839961203
$ws.Range("G11").Value = 0.1740740740740741
$ws.Range("J11").Value = 0.1111111111111111
$ws.Range("K11").Value = 0.2111111111111111
$ws.Range("L11").Value = 0.4851851851851852
$ws.Range("S11").Value = 0.01851851851851852
$ws.Range("G12").Value = 0.6985294117647058
$ws.Range("J12").Value = 0.2205882352941176
$ws.Range("K12").Value = 0.01470588235294118
$ws.Range("L12").Value = 0.02205882352941177
$ws.Range("S12").Value = 0.04411764705882353
$ws.Range("G13").Value = 0.5
$ws.Range("J13").Value = 0.4230769230769231
$ws.Range("S13").Value = 0.07692307692307693
$ws.Range("F15").Value = 0.02762430939226519
$ws.Range("H15").Value = 0.2044198895027624
$ws.Range("I15").Value = 0.06629834254143646
$ws.Range("J15").Value = 0.3370165745856354
$ws.Range("K15").Value = 0.0718232044198895
$ws.Range("O15").Value = 0.05524861878453038
$ws.Range("S15").Value = 0.2375690607734807
$ws.Range("F16").Value = 0.03289473684210526
$ws.Range("H16").Value = 0.1907894736842105
$ws.Range("I16").Value = 0.09210526315789473
$ws.Range("J16").Value = 0.4013157894736842
$ws.Range("K16").Value = 0.1118421052631579
$ws.Range("M16").Value = 0.0131578947368421
$ws.Range("O16").Value = 0.05263157894736842
$ws.Range("S16").Value = 0.1052631578947368
$ws.Range("F17").Value = 0.02639296187683285
$ws.Range("H17").Value = 0.1671554252199413
$ws.Range("I17").Value = 0.06451612903225806
$ws.Range("J17").Value = 0.4340175953079179
$ws.Range("K17").Value = 0.1055718475073314
$ws.Range("M17").Value = 0.01466275659824047
$ws.Range("O17").Value = 0.06744868035190615
$ws.Range("S17").Value = 0.1202346041055719
$ws.Range("F18").Value = 0.008064516129032258
$ws.Range("H18").Value = 0.1129032258064516
$ws.Range("I18").Value = 0.04838709677419355
$ws.Range("J18").Value = 0.4435483870967742
$ws.Range("K18").Value = 0.1370967741935484
$ws.Range("M18").Value = 0.01612903225806452
$ws.Range("O18").Value = 0.08870967741935484
$ws.Range("S18").Value = 0.1451612903225807
$ws.Range("F19").Value = 0.025
$ws.Range("H19").Value = 0.2080357142857143
$ws.Range("I19").Value = 0.07232142857142858
$ws.Range("J19").Value = 0.35
$ws.Range("K19").Value = 0.1098214285714286
$ws.Range("M19").Value = 0.01696428571428571
$ws.Range("N19").Value = 0.002678571428571429
$ws.Range("O19").Value = 0.06160714285714286
$ws.Range("S19").Value = 0.1535714285714286
